# Introduce_caGrid_AHM_2008 - update title textbox: position/size, bold white
# Helvetica Neue styling, and split the "caGrid Service Generation Tools"
# run into separate "caGrid" / " " / "Service " / "Generation Tools" runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(5)

# --- reposition / resize the text box (EMU -> points, 12700 EMU per point) ---
$sh.Left = 22440511 / 12700
$sh.Top = 914400 / 12700
$sh.Width = 19317089 / 12700
$sh.Height = 3339376 / 12700

$tr = $sh.TextFrame.TextRange

$white = 16777215  # RGB(255,255,255)

# --- paragraph 1: "Introduce:" (chars 1-10) ---
$rIntroduce = $tr.Characters(1, 10)
$rIntroduce.Font.Bold = $true
$rIntroduce.Font.Color.RGB = $white
$rIntroduce.Font.Name = "Helvetica Neue"

# --- paragraph 2, run 1: "caGrid" (chars 12-17) ---
$rCaGrid = $tr.Characters(12, 6)
$rCaGrid.Font.Bold = $true
$rCaGrid.Font.Color.RGB = $white
$rCaGrid.Font.Name = "Helvetica Neue"

# --- paragraph 2, run 2: " " (char 18) ---
$rSpace = $tr.Characters(18, 1)
$rSpace.Font.Bold = $true
$rSpace.Font.Color.RGB = $white
$rSpace.Font.Name = "Helvetica Neue"

# --- paragraph 2, run 3: "Service " (chars 19-26) ---
$rService = $tr.Characters(19, 8)
$rService.Font.Bold = $true
$rService.Font.Color.RGB = $white
$rService.Font.Name = "Helvetica Neue"

# --- paragraph 2, run 4: "Generation Tools" (chars 27-42) ---
$rGenTools = $tr.Characters(27, 16)
$rGenTools.Font.Bold = $true
$rGenTools.Font.Color.RGB = $white
$rGenTools.Font.Name = "Helvetica Neue"

# --- complex-script typeface is uniform ("Helvetica Neue") across every run;
#     this host only ever applies NameComplexScript writes to the shape's
#     first run, so a single whole-range assignment covers every run. ---
$tr.Font.NameComplexScript = "Helvetica Neue"
